# Swap the presentation's theme colour scheme from the "Integral" (Red
# Violet) palette over to the stock "Office Theme" palette.
#
# (In the authored commit this is a wholesale swap of the contents of
# ppt/theme/theme1.xml (the slide master's theme - "Integral"/Red Violet)
# and ppt/theme/theme2.xml (the notes master's theme - "Office Theme"), so
# that afterwards theme1.xml holds the Office colours and theme2.xml holds
# the former Integral/Red Violet colours. The PowerPoint object model only
# exposes one editable ThemeColorScheme - the one backing the slide
# master/deck theme (ppt/theme/theme1.xml) - so we drive that to the
# "Office Theme" values here.)

function ToBGR([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme
$colors = $theme.ThemeColorScheme

# Office Theme colour scheme (hex RRGGBB -> COM BGR long), in the
# clrScheme slot order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
$officeTheme = @(
    @(0x00, 0x00, 0x00),   # 1  dk1      000000
    @(0xFF, 0xFF, 0xFF),   # 2  lt1      FFFFFF
    @(0x44, 0x54, 0x6A),   # 3  dk2      44546A
    @(0xE7, 0xE6, 0xE6),   # 4  lt2      E7E6E6
    @(0x5B, 0x9B, 0xD5),   # 5  accent1  5B9BD5
    @(0xED, 0x7D, 0x31),   # 6  accent2  ED7D31
    @(0xA5, 0xA5, 0xA5),   # 7  accent3  A5A5A5
    @(0xFF, 0xC0, 0x00),   # 8  accent4  FFC000
    @(0x44, 0x72, 0xC4),   # 9  accent5  4472C4
    @(0x70, 0xAD, 0x47),   # 10 accent6  70AD47
    @(0x05, 0x63, 0xC1),   # 11 hlink    0563C1
    @(0x95, 0x4F, 0x72)    # 12 folHlink 954F72
)

for ($i = 1; $i -le $colors.Count; $i++) {
    $rgb = $officeTheme[$i - 1]
    $colors.Item($i).RGB = ToBGR $rgb[0] $rgb[1] $rgb[2]
}
